$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.178.22"
$ws.Range("E2").Value = "  -2.84%  "

$ws.Range("D3").Value = "2.368.57"
$ws.Range("E3").Value = "  -3.59%  "

$ws.Range("E4").Value = "  -0.35%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "501.60"
$ws.Range("E5").Value = "  -1.81%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "128.98"
$ws.Range("E6").Value = "  -3.47%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.08%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.545"
$ws.Range("E8").Value = "  -2.28%  "

$ws.Range("D9").Value = "2.373.66"
$ws.Range("E9").Value = "  -3.41%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0983"
$ws.Range("E10").Value = "  +0.53%  "

$ws.Range("E11").Value = "  +0.26%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.84"
$ws.Range("E12").Value = "  +4.75%  "

$ws.Range("E13").Value = "  -0.03%  "

$ws.Range("D14").Value = "2.787.59"
$ws.Range("E14").Value = "  -3.56%  "

$ws.Range("D15").Value = "56.084.77"
$ws.Range("E15").Value = "  -2.93%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.47"
$ws.Range("E16").Value = "  -1.87%  "

$ws.Range("E17").Value = "  -1.20%  "

$ws.Range("D18").Value = "2.409.57"
$ws.Range("E18").Value = "  -3.33%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.01"
$ws.Range("E19").Value = "  -2.99%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.04"
$ws.Range("E20").Value = "  -2.32%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "306.78"
$ws.Range("E21").Value = "  -2.56%  "

$ws.Range("E22").Value = "  -1.97%  "

$ws.Range("E23").Value = "  -0.17%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.93"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.998"
$ws.Range("E25").Value = "  +0.33%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.369"
$ws.Range("E26").Value = "  -3.36%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.146"
$ws.Range("E27").Value = "  -5.91%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.20"
$ws.Range("E28").Value = "  -4.86%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "171.66"
$ws.Range("E29").Value = "  -0.84%  "

$ws.Range("E30").Value = "  -3.09%  "

$ws.Range("E31").Value = "  -3.52%  "

$ws.Range("E32").Value = "  +0.04%  "

$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.08"
$ws.Range("E33").Value = "  -4.97%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.74"
$ws.Range("E34").Value = "  -6.91%  "

$ws.Range("B35").Value = "FirstDigitalUSD"
$ws.Range("C35").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.997"
$ws.Range("E35").Value = "  -0.01%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.61"
$ws.Range("E36").Value = "  -2.40%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.18"
$ws.Range("E37").Value = "  -5.73%  "

$ws.Range("E38").Value = "  -2.93%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.09"
$ws.Range("E39").Value = "  -1.71%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.792"
$ws.Range("E40").Value = "  -1.96%  "

$ws.Range("E41").Value = "  -5.70%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "129.81"
$ws.Range("E42").Value = "  -5.28%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.36"
$ws.Range("E43").Value = "  -1.23%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.69"
$ws.Range("E44").Value = "  -4.63%  "

$ws.Range("E45").Value = "  -2.05%  "

$ws.Range("E46").Value = "  -1.73%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "240.00"
$ws.Range("E47").Value = "  -6.65%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0480"
$ws.Range("E48").Value = "  -2.68%  "

$ws.Range("E49").Value = "  -3.93%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.07"
$ws.Range("E50").Value = "  -0.79%  "
